$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.371.72'
$ws.Range('E2').Value = '  -1.64%  '
$ws.Range('D3').Value = '2.331.06'
$ws.Range('E3').Value = '  +2.55%  '
$ws.Range('E4').Value = '  +0.31%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '232.76'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +0.54%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '0.647'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +1.95%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '66.55'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  +3.78%  '
$ws.Range('E8').Value = '  +0.09%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.455'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +1.26%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.0969'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -4.08%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '56.64'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -1.32%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '26.88'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -0.54%  '
$ws.Range('D13').Value = '2.680.92'
$ws.Range('E13').Value = '  +2.54%  '
$ws.Range('E14').Value = '  -1.31%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '15.51'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -1.33%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '6.23'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +1.52%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '0.853'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +1.46%  '
$ws.Range('D18').Value = '2.337.73'
$ws.Range('E18').Value = '  +2.76%  '
$ws.Range('D19').Value = '43.292.17'
$ws.Range('E19').Value = '  -1.67%  '
$ws.Range('D20').Value = '0.0₃0981'
$ws.Range('E20').Value = '  -2.81%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '74.12'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +0.32%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '6.25'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +2.18%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '248.95'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  -2.03%  '
$ws.Range('E24').Value = '  +13.42%  '
$ws.Range('E25').Value = '  -0.17%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '2.42'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -1.77%  '
$ws.Range('E27').Value = '  -0.96%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '9.93'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '174.24'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +1.53%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '22.13'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +6.10%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '1.46'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +4.84%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '0.128'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -7.69%  '
$ws.Range('E33').Value = '  +0.38%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '5.01'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +3.90%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '0.0688'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -1.58%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '4.96'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +0.93%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '2.52'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +8.68%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '6.51'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -0.57%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '3.61'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -4.75%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.0253'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -2.61%  '
$ws.Range('E41').Value = '  +9.74%  '
$ws.Range('E42').Value = '  +0.18%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '18.26'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +3.65%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '1.17'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +8.23%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '99.29'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +0.93%  '
$ws.Range('E46').Value = '  -0.09%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.0948'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -3.18%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '4.33'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -2.97%  '
$ws.Range('D49').Value = '1.447.56'
$ws.Range('E49').Value = '  -0.15%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '9.94'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -5.56%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '2.29'
$cell.NumberFormat = 'General'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -0.34%  '
